$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46074 -> 46075) for every data row (rows 2 through 378).
$ws.Range("C2:C378").Value = 46075
